# Fruta / hortaliza, semanal
# Insert a new weekly record as row 404, pushing the existing rows
# 404-423 down to 405-424.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 404..423 down by one (new blank row 404, formatting carried
# from the row that used to be there).
$ws.Rows.Item(404).Insert()

# Populate the new record in row 404.
$ws.Cells.Item(404, 1).Value = 7
$ws.Cells.Item(404, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(404, 3).Value = "Ñuble"
$ws.Cells.Item(404, 4).Value = 45008
$ws.Cells.Item(404, 5).Value = 16
$ws.Cells.Item(404, 6).Value = 100112008
$ws.Cells.Item(404, 7).Value = "Coliflor"
$ws.Cells.Item(404, 8).Value = "Sin especificar"
$ws.Cells.Item(404, 9).Value = "Primera"
$ws.Cells.Item(404, 10).Value = 160
$ws.Cells.Item(404, 11).Value = 1200
$ws.Cells.Item(404, 12).Value = 1300
$ws.Cells.Item(404, 13).Value = 1250
$ws.Cells.Item(404, 14).Value = "`$/unidad"
$ws.Cells.Item(404, 15).Value = "Región del Maule"
$ws.Cells.Item(404, 16).Value = 1250
$ws.Cells.Item(404, 17).Value = 1
$ws.Cells.Item(404, 18).Value = "Hortaliza"
